# Feria Lagunitas de Puerto Montt - Aji
# A new weekly price record was inserted into the dataset at row 192,
# shifting all subsequent rows (old 192-229) down by one (new 193-230).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 192, pushing existing
# rows 192..229 down to 193..230.
$ws.Rows("192:192").Insert()

# Populate the newly inserted row 192 with the new data record.
$ws.Cells.Item(192, 1).Value2  = 4
$ws.Cells.Item(192, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(192, 3).Value2  = "Los Lagos"
$ws.Cells.Item(192, 4).Value2  = 44637
$ws.Cells.Item(192, 5).Value2  = 10
$ws.Cells.Item(192, 6).Value2  = 100112021
$ws.Cells.Item(192, 7).Value2  = "Ají"
$ws.Cells.Item(192, 8).Value2  = "Inferno"
$ws.Cells.Item(192, 9).Value2  = "Primera"
$ws.Cells.Item(192, 10).Value2 = 60
$ws.Cells.Item(192, 11).Value2 = 23000
$ws.Cells.Item(192, 12).Value2 = 23000
$ws.Cells.Item(192, 13).Value2 = 23000
$ws.Cells.Item(192, 14).Value2 = "$/caja 15 kilos"
$ws.Cells.Item(192, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(192, 16).Value2 = 1533
$ws.Cells.Item(192, 17).Value2 = 15
$ws.Cells.Item(192, 18).Value2 = "Hortaliza"
